# Apply the commit's edits to the "Embedded_Qs" workbook.
#
# Summary of the change:
#  - Sheet "3_" (the Yes/No outlier question) is converted into a
#    True/False question, matching the style used elsewhere in the
#    workbook (e.g. the "True_False" sheet): the question statement is
#    duplicated into column A of the answer row, the answer becomes a
#    single "T", and the old "No" row is cleared out.
#  - As a side effect of the no-longer-used "Is this an outlier?", "Yes"
#    and "No" shared strings being dropped, the shared-string indices
#    referenced by sheets "1_" and "2_" shift down - this happens
#    automatically when we rewrite the workbook, since Excel recomputes
#    the shared string table from actual cell contents.
#  - The active sheet moves from "4_" (index 4) to "3_" (index 3).

$wb = $excel.ActiveWorkbook

$question = "Imagine you measure the torque on motor shaft, and find it has a mean value of 123.0 Nm with a standard deviation of 6.0 Nm.   One of the 8 readings has a value of 108.0 Nm.   True or False: This data point is an outlier."

# --- Sheet "3_": rework the Yes/No question into a True/False question ---
$s3_ = $wb.Worksheets.Item("3_")

$s3_.Range("A1").Value = $question
$s3_.Rows.Item(1).RowHeight = 105

$s3_.Range("A2").Value = $question
$s3_.Range("B2").Value = "T"
$s3_.Rows.Item(2).RowHeight = 105
# C2 (the comment) keeps its existing text.

# Old row 3 ("No" / "N") is no longer needed.
$s3_.Range("A3").Value = ""
$s3_.Range("B3").Value = ""

# --- Update selections on the other affected sheets (without leaving ---
# --- them as the active tab) before finally activating sheet "3_".   ---

# Sheet "2_" gains a remembered selection at C15.
$s2_ = $wb.Worksheets.Item("2_")
$null = $s2_.Range("C15").Select()

# Sheet "3_" becomes the active tab, with its cursor left at C16.
$null = $s3_.Range("C16").Select()
